# Refresh the "cryptos" price table (rows 2-51 of the active sheet) with a
# newer scrape, per "Updated cryptos list on Thu Jun  8 19:22:26 UTC 2023
# with GitHub Actions".
#
# Columns:
#   A = rank (unchanged)      C = coin link
#   B = coin name              D = price        E = 1h volume/change
#
# Two rows swapped position in the source feed (Dai/ShibaInu at rows 18-19,
# Stellar/Filecoin at rows 30-31), so B/C are rewritten there too; every
# other touched row only gets new D/E figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E hold plain text (prices like "26.558.46" or "1.000", deltas like
# "  +0.27%  "), not numbers -- several of the new values parse as valid
# numbers, so force Text format before writing, then drop back to the
# workbook's default "Normal" style so the cell ends up exactly as
# unstyled as it started (just like the other inlineStr cells on the sheet).
function Set-TextCell([string]$addr, [string]$value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

Set-TextCell 'D2' '26.558.46'
Set-TextCell 'E2' '  +0.27%  '
Set-TextCell 'D3' '1.847.75'
Set-TextCell 'E3' '  -0.07%  '
Set-TextCell 'E4' '  -0.06%  '
Set-TextCell 'D5' '264.15'
Set-TextCell 'E5' '  +1.41%  '
Set-TextCell 'E6' '  +0.00%  '
Set-TextCell 'D7' '0.5213'
Set-TextCell 'E7' '  +0.99%  '
Set-TextCell 'D8' '0.3233'
Set-TextCell 'E8' '  -0.53%  '
Set-TextCell 'E9' '  +0.63%  '
Set-TextCell 'D10' '18.80'
Set-TextCell 'E10' '  -0.48%  '
Set-TextCell 'D11' '0.7810'
Set-TextCell 'E11' '  +1.28%  '
Set-TextCell 'D12' '0.07767'
Set-TextCell 'E12' '  +0.65%  '
Set-TextCell 'D13' '1.828.38'
Set-TextCell 'E13' '  -3.24%  '
Set-TextCell 'D14' '88.56'
Set-TextCell 'E14' '  -0.04%  '
Set-TextCell 'D15' '5.031'
Set-TextCell 'E15' '  +0.04%  '
Set-TextCell 'E16' '  -0.10%  '
Set-TextCell 'D17' '13.99'
Set-TextCell 'E17' '  -0.75%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D18' '0.000007967'
Set-TextCell 'E18' '  +0.58%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D19' '1.000'
Set-TextCell 'E19' '  +0.03%  '
Set-TextCell 'D20' '26.612.01'
Set-TextCell 'E20' '  +0.28%  '
Set-TextCell 'D21' '4.633'
Set-TextCell 'E21' '  +2.46%  '
Set-TextCell 'D22' '9.467'
Set-TextCell 'E22' '  -0.70%  '
Set-TextCell 'D23' '6.017'
Set-TextCell 'E23' '  +1.66%  '
Set-TextCell 'D24' '143.12'
Set-TextCell 'E24' '  -0.96%  '
Set-TextCell 'D25' '2.174'
Set-TextCell 'D26' '1.684'
Set-TextCell 'E26' '  +2.17%  '
Set-TextCell 'D27' '17.02'
Set-TextCell 'E27' '  +0.27%  '
Set-TextCell 'D28' '111.82'
Set-TextCell 'D29' '4.196'
Set-TextCell 'E29' '  -0.32%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D30' '4.119'
Set-TextCell 'E30' '  -1.30%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D31' '0.08747'
Set-TextCell 'E31' '  -0.10%  '
Set-TextCell 'D32' '0.04845'
Set-TextCell 'E32' '  +0.66%  '
Set-TextCell 'E33' '  +5.07%  '
Set-TextCell 'D34' '1.130'
Set-TextCell 'E34' '  -0.32%  '
Set-TextCell 'D35' '2.861'
Set-TextCell 'E35' '  +0.64%  '
Set-TextCell 'D36' '3.107'
Set-TextCell 'E36' '  -0.36%  '
Set-TextCell 'E37' '  -0.16%  '
Set-TextCell 'D38' '2.221'
Set-TextCell 'E38' '  +0.69%  '
Set-TextCell 'D39' '0.4868'
Set-TextCell 'E39' '  -0.47%  '
Set-TextCell 'D40' '111.26'
Set-TextCell 'E40' '  -1.56%  '
Set-TextCell 'D41' '0.8944'
Set-TextCell 'E41' '  -0.47%  '
Set-TextCell 'D42' '6.035'
Set-TextCell 'E42' '  -1.66%  '
Set-TextCell 'D43' '1.000'
Set-TextCell 'E43' '  +0.04%  '
Set-TextCell 'D44' '7.638'
Set-TextCell 'E44' '  -1.82%  '
Set-TextCell 'D45' '0.4211'
Set-TextCell 'E45' '  -0.24%  '
Set-TextCell 'E46' '  -0.04%  '
Set-TextCell 'D47' '9.028'
Set-TextCell 'E47' '  -1.38%  '
Set-TextCell 'D48' '0.1241'
Set-TextCell 'E48' '  -1.28%  '
Set-TextCell 'D49' '35.03'
Set-TextCell 'E49' '  -0.50%  '
Set-TextCell 'D50' '0.8884'
Set-TextCell 'D51' '59.97'
Set-TextCell 'E51' '  +1.28%  '
